# PROS-13075 - CCRU - POS KPI 2020 change
#
# The KPI list in rows 2-14 is re-sorted alphabetically by KPI name
# (columns B, D, E and F all repeat the same KPI text) and a typo is fixed
# ("and-or" -> "and/or") on the "cooler doors" KPI. The column widths are
# also trimmed down slightly, the row-14 auto height is restored (it no
# longer needs the slightly taller 13.8pt it had before) and the sheet
# selection is moved to the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New alphabetically-sorted KPI names for rows 2-14 (columns B, D, E, F all
# carry the same text; column A keeps its "Benchmark 2020" label and column
# C keeps its weight of 1, both untouched).
$names = @(
    "CCH coolers quality",
    "CCH coolers quality (Prime Pos/Max15/Merch STD/Occupancy/Lights&chilled)",
    "CCH products present in Customers menu",
    "CCH shelf share in Energy",
    "CCH shelf share in Juice",
    "CCH shelf share in SSD",
    "CCH shelf share in Tea",
    "CCH shelf share in Water",
    "Number of CCH activation points in NARTD",
    "Number of CCH cooler doors and/or equivalent in Customer coolers",
    "Number of CCH displays points of interaction",
    "Number of NCB core assortment available in-store",
    "Number of SSD core assortment available in-store"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $name = $names[$i]
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 4).Value = $name
    $ws.Cells.Item($row, 5).Value = $name
    $ws.Cells.Item($row, 6).Value = $name
}

# Column width adjustments (slightly narrower than before).
$ws.Columns.Item(1).ColumnWidth = 16.6801619433198
$ws.Columns.Item(2).ColumnWidth = 61
$ws.Columns.Item(3).ColumnWidth = 17.7813765182186
$ws.Columns.Item(4).ColumnWidth = 61
$ws.Columns.Item(5).ColumnWidth = 61
$ws.Columns.Item(6).ColumnWidth = 61

# Row 14 no longer needs a custom (taller) height - let it size back to the
# sheet's default row height.
$ws.Rows.Item(14).AutoFit()

# Select the whole table.
$ws.Range("A1:F14").Select()
